# database connecting, but not returning data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "done" column (F) as TRUE for rows 3-5, matching the existing
# style used by F2 (centered, no special fill).
$ws.Range("F3:F5").Value = $true
$ws.Range("F3:F5").HorizontalAlignment = -4108  # xlCenter

# reimbursements:post (row 8) is now ASSUME for every role, highlighted red.
$ws.Range("B8:D8").Value = "ASSUME"
$ws.Range("B8:D8").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B8:D8").Interior.Color = 255  # RGB(255,0,0) red

# Move the active selection to A8.
$ws.Range("A8").Select()
